$d = $word.ActiveDocument

function FindRange($searchText) {
    $r = $d.Content
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: $searchText"
    }
    return $r
}

# 1) Bold the first paragraph ("Hooray, you're here! ...")
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Font.Bold = 1

# 2) Color "CTRL+K = Translated" blue (4175D4) in the shortcut-keys paragraph
$r2 = FindRange("CTRL+K = Translated")
$r2.Font.Color = 13923649

# 3) Insert the (moved) _GoBack bookmark between "...pressing C" and "TRL+H!"
$r3 = FindRange("You can see other shortcuts by pressing C")
$r3.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r3)

# 4) Replace "{placeholders} used" with "formats inside the Word Document is used",
#    then bold+italic just the "formats inside the Word Document" part.
$r4 = FindRange("{placeholders} used")
$r4.Text = "formats inside the Word Document is used"
$r4b = FindRange("formats inside the Word Document")
$r4b.Font.Bold = 1
$r4b.Font.Italic = 1

# 5) Add a comment around "comments" in the developers paragraph, and change
#    "They are imported as" -> "You can also add"
$r5 = FindRange("comments")
$cmt = $d.Comments.Add($r5, "This is a comment")
$cmt.Author = "Microsoft Office-Anwender"
$cmt.Initial = "Office"
$r5b = FindRange("They are imported as")
$r5b.Text = "You can also add"

# 6) Color the heart character dark red (C00000)
$r6 = FindRange([string][char]0x2764)
$r6.Font.Color = 192
